$d = $word.ActiveDocument

# Locate the "[In Days]" placeholder in the Size cell and replace it with "6".
$findRange = $d.Content
$findRange.Find.Execute("[In Days]", $false, $false, $false, $false, $false,
                         $true, 1, $false, "6", 2)

# Figure out exactly where the replacement text ended up so the _GoBack
# bookmark (which Word drops at the site of the most recent edit) can be
# re-created at the right spot.
$locate = $d.Content
$locate.Find.Execute("6", $false, $false, $false, $false, $false,
                      $true, 1, $false, "", 0)
$newSpot = $d.Range($locate.End, $locate.End)

# Remove the bookmark from its old location (after "Med") if possible, then
# (re)create it right after the newly entered "6".
try {
    $old = $d.Bookmarks("_GoBack")
    $old.Delete()
} catch {
}

$d.Bookmarks.Add("_GoBack", $newSpot)
